$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.555.31"
$ws.Range("E2").Value = "  +2.03%  "

# Row 3
$ws.Range("D3").Value = "1.670.36"
$ws.Range("E3").Value = "  +2.24%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'219.11"
$ws.Range("E5").Value = "  +2.04%  "

# Row 6
$ws.Range("D6").Value = "'0.527"
$ws.Range("E6").Value = "  +1.77%  "

# Row 7
$ws.Range("E7").Value = "  +0.28%  "

# Row 8
$ws.Range("D8").Value = "'29.07"
$ws.Range("E8").Value = "  +1.03%  "

# Row 9
$ws.Range("E9").Value = "  +1.84%  "

# Row 10
$ws.Range("E10").Value = "  +4.69%  "

# Row 11
$ws.Range("E11").Value = "  -0.11%  "

# Row 12
$ws.Range("D12").Value = "1.915.27"
$ws.Range("E12").Value = "  +2.53%  "

# Row 13
$ws.Range("D13").Value = "1.670.44"
$ws.Range("E13").Value = "  +2.07%  "

# Row 14
$ws.Range("E14").Value = "  +7.06%  "

# Row 15
$ws.Range("D15").Value = "'9.93"
$ws.Range("E15").Value = "  +6.81%  "

# Row 16
$ws.Range("E16").Value = "  +3.89%  "

# Row 17
$ws.Range("D17").Value = "30.540.18"
$ws.Range("E17").Value = "  +1.92%  "

# Row 18
$ws.Range("D18").Value = "'65.88"
$ws.Range("E18").Value = "  +2.60%  "

# Row 19
$ws.Range("D19").Value = "'241.71"
$ws.Range("E19").Value = "  +0.15%  "

# Row 20
$ws.Range("E20").Value = "  +2.26%  "

# Row 21
$ws.Range("E21").Value = "  +0.20%  "

# Row 22
$ws.Range("D22").Value = "'4.21"
$ws.Range("E22").Value = "  +1.88%  "

# Row 23
$ws.Range("D23").Value = "'9.93"
$ws.Range("E23").Value = "  +1.19%  "

# Row 24
$ws.Range("E24").Value = "  -0.44%  "

# Row 25
$ws.Range("D25").Value = "'158.81"
$ws.Range("E25").Value = "  +0.57%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'15.74"
$ws.Range("E26").Value = "  +1.61%  "

# Row 27
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.112"
$ws.Range("E27").Value = "  +2.11%  "

# Row 28
$ws.Range("D28").Value = "'6.64"
$ws.Range("E28").Value = "  +1.24%  "

# Row 29
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("D30").Value = "'0.0491"
$ws.Range("E30").Value = "  +0.13%  "

# Row 31
$ws.Range("E31").Value = "  +3.17%  "

# Row 32
$ws.Range("D32").Value = "'3.44"
$ws.Range("E32").Value = "  +1.62%  "

# Row 33
$ws.Range("D33").Value = "'3.28"
$ws.Range("E33").Value = "  +3.25%  "

# Row 34
$ws.Range("D34").Value = "1.494.46"
$ws.Range("E34").Value = "  +4.53%  "

# Row 35
$ws.Range("E35").Value = "  +5.90%  "

# Row 36
$ws.Range("E36").Value = "  -0.99%  "

# Row 37
$ws.Range("D37").Value = "'82.99"
$ws.Range("E37").Value = "  +9.72%  "

# Row 38
$ws.Range("D38").Value = "'0.593"
$ws.Range("E38").Value = "  +7.14%  "

# Row 39
$ws.Range("E39").Value = "  +3.73%  "

# Row 40
$ws.Range("E40").Value = "  -3.29%  "

# Row 41
$ws.Range("E41").Value = "  +0.75%  "

# Row 42
$ws.Range("D42").Value = "'2.00"
$ws.Range("E42").Value = "  +0.75%  "

# Row 43
$ws.Range("E43").Value = "  +0.59%  "

# Row 44
$ws.Range("D44").Value = "'0.0498"
$ws.Range("E44").Value = "  +0.33%  "

# Row 45
$ws.Range("E45").Value = "  +1.74%  "

# Row 46
$ws.Range("E46").Value = "  +0.28%  "

# Row 47
$ws.Range("E47").Value = "  +3.15%  "

# Row 48
$ws.Range("D48").Value = "1.808.40"
$ws.Range("E48").Value = "  +1.88%  "

# Row 49
$ws.Range("D49").Value = "'49.30"
$ws.Range("E49").Value = "  -3.93%  "

# Row 50
$ws.Range("D50").Value = "'93.52"
$ws.Range("E50").Value = "  +3.43%  "

# Row 51
$ws.Range("E51").Value = "  -1.72%  "
